# Update column F ("dSF") values on Sheet1 to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -4
    4  = -5
    6  = -7
    7  = -2
    9  = -1
    10 = 2
    14 = -6
    16 = -2
    20 = -6
    27 = -8
    28 = -7
    29 = -2
    32 = -5
    34 = 11
    37 = -4
    38 = -3
    41 = 0
    43 = 0
    45 = 0
    47 = 3
    54 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
